$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.07443033333333333
$ws.Range("H2").Value = 0.223291
$ws.Range("I2").Value = 0.07586947613633815
$ws.Range("J2").Value = 0.07586947613633817
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03442933333333333
$ws.Range("N2").Value = 0.103288
$ws.Range("O2").Value = 0.13994583074207
$ws.Range("P2").Value = 0.1399458307420701
$ws.Range("Q2").Value = 0.002562586756444444
$ws.Range("R2").Value = 0.023063280808
$ws.Range("S2").Value = 0.0106176168658655
$ws.Range("T2").Value = 0.01061761686586551
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.07443033333333333
$ws.Range("H3").Value = 0.223291
$ws.Range("I3").Value = 0.07586947613633815
$ws.Range("J3").Value = 0.07586947613633817
$ws.Range("O3").Value = 0.6285530792337177
$ws.Range("P3").Value = 0.6285530792337178
$ws.Range("Q3").Value = 0.01150960902533333
$ws.Range("R3").Value = 0.103586481228
$ws.Range("S3").Value = 0.04768799284534441
$ws.Range("T3").Value = 0.04768799284534443
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.07443033333333333
$ws.Range("H4").Value = 0.223291
$ws.Range("I4").Value = 0.07586947613633815
$ws.Range("J4").Value = 0.07586947613633817
$ws.Range("M4").Value = 0.05695366666666666
$ws.Range("N4").Value = 0.170861
$ws.Range("O4").Value = 0.2315010900242122
$ws.Range("P4").Value = 0.2315010900242122
$ws.Range("Q4").Value = 0.004239080394555555
$ws.Range("R4").Value = 0.038151723551
$ws.Range("S4").Value = 0.01756386642512824
$ws.Range("T4").Value = 0.01756386642512824
$ws.Range("I5").Value = 0.7501574873245639
$ws.Range("J5").Value = 0.7501574873245638
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03442933333333333
$ws.Range("N5").Value = 0.103288
$ws.Range("O5").Value = 0.13994583074207
$ws.Range("P5").Value = 0.1399458307420701
$ws.Range("Q5").Value = 0.02533751042133333
$ws.Range("R5").Value = 0.228037593792
$ws.Range("S5").Value = 0.10498141275102
$ws.Range("T5").Value = 0.10498141275102
$ws.Range("I6").Value = 0.7501574873245639
$ws.Range("J6").Value = 0.7501574873245638
$ws.Range("O6").Value = 0.6285530792337177
$ws.Range("P6").Value = 0.6285530792337178
$ws.Range("S6").Value = 0.4715137985680832
$ws.Range("T6").Value = 0.4715137985680832
$ws.Range("I7").Value = 0.7501574873245639
$ws.Range("J7").Value = 0.7501574873245638
$ws.Range("M7").Value = 0.05695366666666666
$ws.Range("N7").Value = 0.170861
$ws.Range("O7").Value = 0.2315010900242122
$ws.Range("P7").Value = 0.2315010900242122
$ws.Range("Q7").Value = 0.04191379800266665
$ws.Range("R7").Value = 0.3772241820239999
$ws.Range("S7").Value = 0.1736622760054607
$ws.Range("T7").Value = 0.1736622760054607
$ws.Range("G8").Value = 0.170673
$ws.Range("H8").Value = 0.512019
$ws.Range("I8").Value = 0.173973036539098
$ws.Range("J8").Value = 0.173973036539098
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03442933333333333
$ws.Range("N8").Value = 0.103288
$ws.Range("O8").Value = 0.13994583074207
$ws.Range("P8").Value = 0.1399458307420701
$ws.Range("Q8").Value = 0.005876157607999999
$ws.Range("R8").Value = 0.05288541847200001
$ws.Range("S8").Value = 0.02434680112518457
$ws.Range("T8").Value = 0.02434680112518458
$ws.Range("G9").Value = 0.170673
$ws.Range("H9").Value = 0.512019
$ws.Range("I9").Value = 0.173973036539098
$ws.Range("J9").Value = 0.173973036539098
$ws.Range("O9").Value = 0.6285530792337177
$ws.Range("P9").Value = 0.6285530792337178
$ws.Range("Q9").Value = 0.026392190028
$ws.Range("R9").Value = 0.237529710252
$ws.Range("S9").Value = 0.1093512878202901
$ws.Range("T9").Value = 0.1093512878202901
$ws.Range("G10").Value = 0.170673
$ws.Range("H10").Value = 0.512019
$ws.Range("I10").Value = 0.173973036539098
$ws.Range("J10").Value = 0.173973036539098
$ws.Range("M10").Value = 0.05695366666666666
$ws.Range("N10").Value = 0.170861
$ws.Range("O10").Value = 0.2315010900242122
$ws.Range("P10").Value = 0.2315010900242122
$ws.Range("Q10").Value = 0.009720453150999998
$ws.Range("R10").Value = 0.08748407835899999
$ws.Range("S10").Value = 0.04027494759362328
$ws.Range("T10").Value = 0.04027494759362328
